$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to text format so numeric-looking values
# (e.g. "0.9985", "30.297.52", "0.06430") are stored as literal text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$priceUpdates = @{
    2 = "30.297.52"
    3 = "1.860.40"
    4 = "0.9985"
    5 = "233.45"
    6 = "0.9988"
    7 = "0.4723"
    8 = "0.2742"
    9 = "0.06430"
    10 = "1.828.48"
    11 = "0.07439"
    12 = "16.27"
    13 = "5.011"
    14 = "85.60"
    15 = "0.6341"
    16 = "30.251.82"
    17 = "0.9994"
    18 = "233.21"
    19 = "12.76"
    20 = "0.000007382"
    21 = "2.093.22"
    22 = "0.9991"
    23 = "5.042"
    24 = "6.012"
    25 = "9.269"
    26 = "165.03"
    27 = "17.94"
    28 = "1.893"
    29 = "0.1032"
    30 = "1.379"
    31 = "4.145"
    32 = "3.940"
    33 = "0.04897"
    34 = "1.148"
    35 = "0.7253"
    36 = "0.9990"
    37 = "2.691"
    38 = "0.01910"
    40 = "0.9059"
    41 = "1.975"
    42 = "105.73"
    43 = "0.9986"
    44 = "0.4119"
    45 = "5.527"
    46 = "7.181"
    47 = "61.16"
    48 = "0.1200"
    49 = "8.710"
    50 = "1.407"
    51 = "33.32"
}

$volumeUpdates = @{
    2 = "  -1.64%  "
    3 = "  -1.21%  "
    4 = "  -0.13%  "
    5 = "  -2.52%  "
    6 = "  -0.08%  "
    7 = "  -2.05%  "
    8 = "  -3.76%  "
    9 = "  -1.73%  "
    10 = "  -7.87%  "
    11 = "  -0.79%  "
    12 = "  -2.30%  "
    13 = "  -2.07%  "
    14 = "  -3.61%  "
    15 = "  -5.28%  "
    16 = "  -1.66%  "
    17 = "  +0.04%  "
    18 = "  +0.58%  "
    19 = "  -4.46%  "
    20 = "  -3.08%  "
    21 = "  -5.31%  "
    22 = "  -0.08%  "
    23 = "  -5.04%  "
    24 = "  -2.93%  "
    25 = "  -0.76%  "
    26 = "  -1.98%  "
    27 = "  -4.21%  "
    28 = "  -2.78%  "
    29 = "  +6.15%  "
    30 = "  -2.71%  "
    31 = "  -5.17%  "
    32 = "  -2.44%  "
    33 = "  -3.41%  "
    34 = "  -5.48%  "
    35 = "  -3.88%  "
    36 = "  -0.60%  "
    37 = "  -0.56%  "
    38 = "  +1.76%  "
    39 = "  -0.46%  "
    40 = "  -1.14%  "
    41 = "  -5.78%  "
    42 = "  -0.86%  "
    43 = "  -0.20%  "
    44 = "  -4.17%  "
    45 = "  -4.81%  "
    46 = "  -2.92%  "
    47 = "  -5.63%  "
    48 = "  -6.86%  "
    49 = "  -2.49%  "
    50 = "  -5.27%  "
    51 = "  -1.78%  "
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $volumeUpdates[$row]
}

Write-Output "Updated $($priceUpdates.Count) price cells and $($volumeUpdates.Count) volume cells."
